# Applies the scheduled-runner market data refresh to the Cuchulainn Profits
# workbook: updates computed price/profit columns (H, I, J, K, L, M, N) for the
# affected leve rows across all crafting-class sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1357.0714
$ws.Range("I19").Value = 1867.2222
$ws.Range("J19").Value = 438.8
$ws.Range("K19").Value = 1867.2222
$ws.Range("L19").Value = 438.8
$ws.Range("M19").Value = -1692.2222
$ws.Range("N19").Value = -788.8
$ws.Range("H28").Value = 1348
$ws.Range("I28").Value = 748
$ws.Range("J28").Value = 3148
$ws.Range("K28").Value = 748
$ws.Range("L28").Value = 3148
$ws.Range("M28").Value = -263
$ws.Range("N28").Value = -4118
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = $null
$ws.Range("N41").Value = $null
$ws.Range("H62").Value = 6522.4546
$ws.Range("I62").Value = 5729.8
$ws.Range("K62").Value = 5729.8
$ws.Range("M62").Value = -5105.8
$ws.Range("H65").Value = 6522.4546
$ws.Range("I65").Value = 5729.8
$ws.Range("K65").Value = 28649
$ws.Range("M65").Value = -25529
$ws.Range("H113").Value = 3997.3333
$ws.Range("I113").Value = 3997.3333
$ws.Range("K113").Value = 3997.3333
$ws.Range("M113").Value = -743.3332999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4398.8
$ws.Range("I45").Value = 4398.8
$ws.Range("K45").Value = 4398.8
$ws.Range("M45").Value = -4021.8
$ws.Range("H61").Value = 9998.666999999999
$ws.Range("J61").Value = 19996
$ws.Range("L61").Value = 19996
$ws.Range("N61").Value = -20420
$ws.Range("H63").Value = 6024
$ws.Range("I63").Value = 6024
$ws.Range("K63").Value = 6024
$ws.Range("M63").Value = -5338
$ws.Range("H66").Value = 6024
$ws.Range("I66").Value = 6024
$ws.Range("K66").Value = 30120
$ws.Range("M66").Value = -26688
$ws.Range("H97").Value = 218.75
$ws.Range("I97").Value = 218.75
$ws.Range("K97").Value = 218.75
$ws.Range("M97").Value = 277.25
$ws.Range("H98").Value = 56977.5
$ws.Range("J98").Value = 56977.5
$ws.Range("L98").Value = 56977.5
$ws.Range("N98").Value = -62967.5
$ws.Range("H132").Value = 11496.667
$ws.Range("I132").Value = 9745
$ws.Range("J132").Value = 15000
$ws.Range("K132").Value = 29235
$ws.Range("L132").Value = 45000
$ws.Range("M132").Value = -26705
$ws.Range("N132").Value = -50060
$ws.Range("H136").Value = 9998.666999999999
$ws.Range("J136").Value = 19996
$ws.Range("L136").Value = 59988
$ws.Range("N136").Value = -65088

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 13111
$ws.Range("I82").Value = 6438.75
$ws.Range("K82").Value = 6438.75
$ws.Range("M82").Value = -6055.75
$ws.Range("H85").Value = 13111
$ws.Range("I85").Value = 6438.75
$ws.Range("K85").Value = 6438.75
$ws.Range("M85").Value = -5112.75
$ws.Range("H94").Value = 555.125
$ws.Range("I94").Value = 555.125
$ws.Range("K94").Value = 555.125
$ws.Range("M94").Value = -104.125
$ws.Range("H99").Value = 958.1818
$ws.Range("I99").Value = 1013.2
$ws.Range("J99").Value = 408
$ws.Range("K99").Value = 1013.2
$ws.Range("L99").Value = 408
$ws.Range("M99").Value = 484.8
$ws.Range("N99").Value = -3404
$ws.Range("H103").Value = 20026.572
$ws.Range("J103").Value = 20026.572
$ws.Range("L103").Value = 20026.572
$ws.Range("N103").Value = -22370.572
$ws.Range("H134").Value = 4679.6
$ws.Range("I134").Value = 4679.6
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 14038.8
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -11503.8
$ws.Range("N134").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 473.1111
$ws.Range("I10").Value = 191.66667
$ws.Range("J10").Value = 1036
$ws.Range("K10").Value = 191.66667
$ws.Range("L10").Value = 1036
$ws.Range("M10").Value = -52.66667000000001
$ws.Range("N10").Value = -1314

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J86").Value = 500
$ws.Range("L86").Value = 1500
$ws.Range("N86").Value = -3872
$ws.Range("J89").Value = 500
$ws.Range("L89").Value = 4500
$ws.Range("N89").Value = -16356
$ws.Range("H131").Value = 910
$ws.Range("I131").Value = 910
$ws.Range("K131").Value = 2730
$ws.Range("M131").Value = 2310

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = $null
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").Value = $null
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").Value = $null
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = $null
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = $null
$ws.Range("H102").Value = 4349.8887
$ws.Range("I102").Value = 3339.2856
$ws.Range("K102").Value = 3339.2856
$ws.Range("M102").Value = -1717.2856
$ws.Range("H113").Value = 3033.8572
$ws.Range("I113").Value = 3083
$ws.Range("J113").Value = 2739
$ws.Range("K113").Value = 3083
$ws.Range("L113").Value = 2739
$ws.Range("M113").Value = -913
$ws.Range("N113").Value = -7079

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2500
$ws.Range("I61").Value = 2500
$ws.Range("K61").Value = 2500
$ws.Range("M61").Value = -2298
$ws.Range("H105").Value = 23016.428
$ws.Range("J105").Value = 23016.428
$ws.Range("L105").Value = 23016.428
$ws.Range("N105").Value = -30004.428
$ws.Range("H113").Value = 2500
$ws.Range("I113").Value = 2500
$ws.Range("K113").Value = 2500
$ws.Range("M113").Value = -330
$ws.Range("H132").Value = 20004
$ws.Range("I132").Value = 20004
$ws.Range("K132").Value = 60012
$ws.Range("M132").Value = -57482

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2995
$ws.Range("I132").Value = 2995
$ws.Range("K132").Value = 2995
$ws.Range("M132").Value = -6455
